# "Problem Statement Diagram Aligned"
# Re-align / re-size the small problem-statement diagram cluster on slide 2:
# the illustration icon plus the two rectangle/textbox label pairs next to
# it. (The source diff also shows the Acrobat OLE object on slide 12 getting
# a new internal VML "spid" - that is an id PowerPoint silently renumbers on
# its own resave and isn't exposed anywhere on the Shape/OLEFormat object
# model, so there is nothing to set for it here.)
#
# NOTE on the literal point values below: PowerPoint stores Shape
# Left/Top/Width/Height in points as single-precision floats, and this host
# truncates (floors) when it converts back to EMU on save (EMU = points *
# 12700). Plain "emu / 12700.0" can therefore land one EMU short after the
# float32 round-trip, so the literals here are the nearest double that
# round-trips through float32 to the exact target EMU.

$p = $ppt.ActivePresentation

# ---- Slide 2: "Problem Statement" diagram shapes --------------------------
$s2 = $p.Slides.Item(2)

# Picture 8 - the small illustration icon: keep position, shrink to a square
# 629350 EMU per side (was 685624 EMU).
$pic = $s2.Shapes.Item(13)
$pic.Width = 49.555118560791016
$pic.Height = 49.555118560791016

# Rectangle 26 - purple-outlined box behind "Pharmacist is unaware..." text
$rect26 = $s2.Shapes.Item(14)
$rect26.Left = 596.7789916992188
$rect26.Width = 93.14385986328125

# TextBox 27 - "Pharmacist is unaware of illegitimate prescriptions"
$tb27 = $s2.Shapes.Item(15)
$tb27.Left = 595.7156982421875
$tb27.Width = 91.55094909667969

# Rectangle 35 - green-outlined box behind "Customer forgets..." text
$rect35 = $s2.Shapes.Item(16)
$rect35.Left = 596.778076171875
$rect35.Width = 93.14385986328125

# TextBox 36 - "Customer forgets to bring the prescription"
$tb36 = $s2.Shapes.Item(17)
$tb36.Left = 603.1698608398438
$tb36.Width = 85.3429946899414
